$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a date serial that was bumped by one day
# (2023-10-06 -> 2023-10-07, serial 45205 -> 45206) for every data row
# (rows 2 through 536).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 536 }

$ws.Range("C2:C$lastRow").Value = 45206
